$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CONTAFEVER N 200MG/5ML SUSP. 120ML (row 11):
#  - current balance (الرصيد الحالي) 9:0 -> 8:0
#  - selling price (سعر البيع) 33.0000 -> 66.0000
#  - transactions count (عدد التعاملات) 1:0 -> 2:0
$ws.Range("H11").Value = "8:0"
$ws.Range("Q11").Value = "2:0"

# P11 holds a numeric-looking string ("66.0000") that must stay text, not be
# auto-coerced into the number 66. Temporarily switch to a text format while
# writing the value, then restore the original "0.00" numeric format so the
# cell's style/format is left exactly as it was.
$p11 = $ws.Range("P11")
$p11.NumberFormat = "@"
$p11.Value = "66.0000"
$p11.NumberFormat = "0.00"

# Grand total of the selling-price column reflects the +33.0000 change
$ws.Range("P27").Value = 766.88

# Footer timestamp updated to the new export time
$ws.Range("A28").Value = "Thursday, 9 October, 2025 2:47 PM"
